$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fgf13"
$ws.Range("C2").Value = "Scn8a"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.14444
$ws.Range("H2").Value = 0.43332
$ws.Range("I2").Value = 0.06801140868936309
$ws.Range("J2").Value = 0.06801140868936309
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.07718933333333333
$ws.Range("N2").Value = 0.231568
$ws.Range("O2").Value = 0.06450640626545157
$ws.Range("P2").Value = 0.06450640626545157
$ws.Range("Q2").Value = 0.01114922730666667
$ws.Range("R2").Value = 0.10034304576
$ws.Range("S2").Value = 0.004387171559601718
$ws.Range("T2").Value = 0.004387171559601718

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fgf13"
$ws.Range("C3").Value = "Scn8a"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.14444
$ws.Range("H3").Value = 0.43332
$ws.Range("I3").Value = 0.06801140868936309
$ws.Range("J3").Value = 0.06801140868936309
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.654434
$ws.Range("N3").Value = 1.963302
$ws.Range("O3").Value = 0.546904392808046
$ws.Range("P3").Value = 0.546904392808046
$ws.Range("Q3").Value = 0.09452644696
$ws.Range("R3").Value = 0.85073802264
$ws.Range("S3").Value = 0.03719573817327599
$ws.Range("T3").Value = 0.03719573817327599

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fgf13"
$ws.Range("C4").Value = "Scn8a"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.14444
$ws.Range("H4").Value = 0.43332
$ws.Range("I4").Value = 0.06801140868936309
$ws.Range("J4").Value = 0.06801140868936309
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4649916666666667
$ws.Range("N4").Value = 1.394975
$ws.Range("O4").Value = 0.3885892009265024
$ws.Range("P4").Value = 0.3885892009265024
$ws.Range("Q4").Value = 0.06716339633333335
$ws.Range("R4").Value = 0.6044705670000001
$ws.Range("S4").Value = 0.02642849895648539
$ws.Range("T4").Value = 0.02642849895648539

$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Fgf13"
$ws.Range("C5").Value = "Scn8a"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.979321333333333
$ws.Range("H5").Value = 5.937964
$ws.Range("I5").Value = 0.9319885913106368
$ws.Range("J5").Value = 0.9319885913106369
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.07718933333333333
$ws.Range("N5").Value = 0.231568
$ws.Range("O5").Value = 0.06450640626545157
$ws.Range("P5").Value = 0.06450640626545157
$ws.Range("Q5").Value = 0.1527824941724444
$ws.Range("R5").Value = 1.375042447552
$ws.Range("S5").Value = 0.06011923470584984
$ws.Range("T5").Value = 0.06011923470584985

$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Fgf13"
$ws.Range("C6").Value = "Scn8a"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.979321333333333
$ws.Range("H6").Value = 5.937964
$ws.Range("I6").Value = 0.9319885913106368
$ws.Range("J6").Value = 0.9319885913106369
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.654434
$ws.Range("N6").Value = 1.963302
$ws.Range("O6").Value = 0.546904392808046
$ws.Range("P6").Value = 0.546904392808046
$ws.Range("Q6").Value = 1.295335177458667
$ws.Range("R6").Value = 11.658016597128
$ws.Range("S6").Value = 0.50970865463477
$ws.Range("T6").Value = 0.5097086546347701

$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Fgf13"
$ws.Range("C7").Value = "Scn8a"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.979321333333333
$ws.Range("H7").Value = 5.937964
$ws.Range("I7").Value = 0.9319885913106368
$ws.Range("J7").Value = 0.9319885913106369
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.4649916666666667
$ws.Range("N7").Value = 1.394975
$ws.Range("O7").Value = 0.3885892009265024
$ws.Range("P7").Value = 0.3885892009265024
$ws.Range("Q7").Value = 0.9203679256555556
$ws.Range("R7").Value = 8.2833113309
$ws.Range("S7").Value = 0.362160701970017
$ws.Range("T7").Value = 0.362160701970017

